# Automatic push - Update excel files
#
# The workbook's "Feuil1" sheet pulls live crypto/asset prices through
# external-workbook references (e.g. =[2]ETH!J4, =-[2]BIGTIME!$C$4, ...).
# The source workbook ("Historique d'achats.xlsx") was refreshed, which
# changed the cached results that those external references resolve to.
# Update the resulting values here so every dependent formula on the sheet
# (percentages, "Others" roll-ups, totals, the pie chart feeding off
# Feuil1!N8:N23, ...) recalculates consistently from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("C12").Value = 2148.2715781179131   # [2]ETH!J4
$ws.Range("C13").Value = 1514.3989186470469   # [2]BTC!J4
$ws.Range("C14").Value = 465.29898996667129   # [2]SOL!J4
$ws.Range("C17").Value = 298.40244046932077   # [2]BNB!J4
$ws.Range("C19").Value = 45.666666666666664   # -[2]BIGTIME!$C$4
$ws.Range("C20").Value = 61.994007660681909   # [2]AVAX!$J$4
$ws.Range("C21").Value = 66.406986848401246   # [2]MATIC!$J$4
$ws.Range("C22").Value = 69.474184778617911   # [2]NEAR!$J$4
$ws.Range("C24").Value = 56.158440351073835   # [2]DOT!$J$4
$ws.Range("C25").Value = 52.33243735765722    # [2]LUNC!J4
$ws.Range("C27").Value = 51.754911112059354   # [2]ADA!$J$4
$ws.Range("C28").Value = 46.104941497289936   # [2]MINA!$J$4
$ws.Range("C29").Value = 21.471098521620156   # [2]TIA!$J$4
$ws.Range("C30").Value = 33.906804324237434   # [2]APE!$J$4
$ws.Range("C31").Value = 22.882289746475138   # [2]DYDX!$J$4
$ws.Range("C32").Value = 18.406194648385924   # [2]UNI!$J$4
$ws.Range("C33").Value = 20.90356376401936    # [2]LDO!$J$4
$ws.Range("C34").Value = 17.379642935991598   # [2]XRP!$J$4
$ws.Range("C35").Value = 15.984983110395214   # [2]SHIB!$J$4
$ws.Range("C36").Value = 14.380564466888945   # [2]ICP!$J$4
$ws.Range("C37").Value = 13.335143111664975   # [2]LINK!$J$4
$ws.Range("C38").Value = 13.50809005440578    # [2]ATOM!$J$4
$ws.Range("C39").Value = 12.225145303102996   # [2]LTC!$J$4
$ws.Range("C40").Value = 13.433690499603543   # [2]ALGO!$J$4
$ws.Range("C42").Value = 6.2190873745437774   # [2]EGLD!$J$4
$ws.Range("C43").Value = 5.9016605154634201   # [2]DOGE!$J$4
$ws.Range("C44").Value = 5.174319884042192    # [2]LUNA!J4
$ws.Range("C45").Value = 7.5611245623061114   # [2]GRT!$J$4
$ws.Range("C46").Value = 2.9141276206016897   # [2]AMP!$J$4
$ws.Range("C47").Value = 3.5921101069419121   # [2]ACE!$J$4
$ws.Range("C48").Value = 3.3691880732865762   # [2]SEI!$J$4
$ws.Range("C49").Value = 3.2504687734656277   # [2]SHPING!$J$4
$ws.Range("C50").Value = 2.4101974840058387   # [2]KAVA!$J$4
$ws.Range("C51").Value = 2.7216952652620194   # [2]POLIS!J4
$ws.Range("C52").Value = 2.0901713503606412   # [2]MEME!$J$4
$ws.Range("C53").Value = 1.3282756559711013   # [2]TRX!$J$4
$ws.Range("C54").Value = 0.59652999419804686  # [2]ATLAS!O47

$excel.CalculateFull()
$wb.Save()
